$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 69.14491266666667
$ws.Cells.Item(2, 8).Value = 207.434738
$ws.Cells.Item(2, 9).Value = 0.8044999916189329
$ws.Cells.Item(2, 10).Value = 0.8044999916189327
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 69.14491266666667
$ws.Cells.Item(2, 14).Value = 207.434738
$ws.Cells.Item(2, 15).Value = 0.8044999916189329
$ws.Cells.Item(2, 16).Value = 0.8044999916189327
$ws.Cells.Item(2, 17).Value = 4781.018947680961
$ws.Cells.Item(2, 18).Value = 43029.17052912865
$ws.Cells.Item(2, 19).Value = 0.6472202365148632
$ws.Cells.Item(2, 20).Value = 0.6472202365148628
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 69.14491266666667
$ws.Cells.Item(3, 8).Value = 207.434738
$ws.Cells.Item(3, 9).Value = 0.8044999916189329
$ws.Cells.Item(3, 10).Value = 0.8044999916189327
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 13.76110233333333
$ws.Cells.Item(3, 14).Value = 41.283307
$ws.Cells.Item(3, 15).Value = 0.1601102132445234
$ws.Cells.Item(3, 16).Value = 0.1601102132445233
$ws.Cells.Item(3, 17).Value = 951.5102190353963
$ws.Cells.Item(3, 18).Value = 8563.591971318567
$ws.Cells.Item(3, 19).Value = 0.1288086652133246
$ws.Cells.Item(3, 20).Value = 0.1288086652133245
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 69.14491266666667
$ws.Cells.Item(4, 8).Value = 207.434738
$ws.Cells.Item(4, 9).Value = 0.8044999916189329
$ws.Cells.Item(4, 10).Value = 0.8044999916189327
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 3.041671
$ws.Cells.Item(4, 14).Value = 9.125013000000001
$ws.Cells.Item(4, 15).Value = 0.03538979513654388
$ws.Cells.Item(4, 16).Value = 0.03538979513654388
$ws.Cells.Item(4, 17).Value = 210.3160756557327
$ws.Cells.Item(4, 18).Value = 1892.844680901594
$ws.Cells.Item(4, 19).Value = 0.02847108989074531
$ws.Cells.Item(4, 20).Value = 0.02847108989074529
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 13.76110233333333
$ws.Cells.Item(5, 8).Value = 41.283307
$ws.Cells.Item(5, 9).Value = 0.1601102132445234
$ws.Cells.Item(5, 10).Value = 0.1601102132445233
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 69.14491266666667
$ws.Cells.Item(5, 14).Value = 207.434738
$ws.Cells.Item(5, 15).Value = 0.8044999916189329
$ws.Cells.Item(5, 16).Value = 0.8044999916189327
$ws.Cells.Item(5, 17).Value = 951.5102190353963
$ws.Cells.Item(5, 18).Value = 8563.591971318567
$ws.Cells.Item(5, 19).Value = 0.1288086652133246
$ws.Cells.Item(5, 20).Value = 0.1288086652133245
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 13.76110233333333
$ws.Cells.Item(6, 8).Value = 41.283307
$ws.Cells.Item(6, 9).Value = 0.1601102132445234
$ws.Cells.Item(6, 10).Value = 0.1601102132445233
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 13.76110233333333
$ws.Cells.Item(6, 14).Value = 41.283307
$ws.Cells.Item(6, 15).Value = 0.1601102132445234
$ws.Cells.Item(6, 16).Value = 0.1601102132445233
$ws.Cells.Item(6, 17).Value = 189.3679374284721
$ws.Cells.Item(6, 18).Value = 1704.311436856249
$ws.Cells.Item(6, 19).Value = 0.02563528038520675
$ws.Cells.Item(6, 20).Value = 0.02563528038520673
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 13.76110233333333
$ws.Cells.Item(7, 8).Value = 41.283307
$ws.Cells.Item(7, 9).Value = 0.1601102132445234
$ws.Cells.Item(7, 10).Value = 0.1601102132445233
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 3.041671
$ws.Cells.Item(7, 14).Value = 9.125013000000001
$ws.Cells.Item(7, 15).Value = 0.03538979513654388
$ws.Cells.Item(7, 16).Value = 0.03538979513654388
$ws.Cells.Item(7, 17).Value = 41.85674589533234
$ws.Cells.Item(7, 18).Value = 376.710713057991
$ws.Cells.Item(7, 19).Value = 0.005666267645992037
$ws.Cells.Item(7, 20).Value = 0.005666267645992034
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3.041671
$ws.Cells.Item(8, 8).Value = 9.125013000000001
$ws.Cells.Item(8, 9).Value = 0.03538979513654388
$ws.Cells.Item(8, 10).Value = 0.03538979513654388
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 69.14491266666667
$ws.Cells.Item(8, 14).Value = 207.434738
$ws.Cells.Item(8, 15).Value = 0.8044999916189329
$ws.Cells.Item(8, 16).Value = 0.8044999916189327
$ws.Cells.Item(8, 17).Value = 210.3160756557327
$ws.Cells.Item(8, 18).Value = 1892.844680901594
$ws.Cells.Item(8, 19).Value = 0.02847108989074531
$ws.Cells.Item(8, 20).Value = 0.02847108989074529
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3.041671
$ws.Cells.Item(9, 8).Value = 9.125013000000001
$ws.Cells.Item(9, 9).Value = 0.03538979513654388
$ws.Cells.Item(9, 10).Value = 0.03538979513654388
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 13.76110233333333
$ws.Cells.Item(9, 14).Value = 41.283307
$ws.Cells.Item(9, 15).Value = 0.1601102132445234
$ws.Cells.Item(9, 16).Value = 0.1601102132445233
$ws.Cells.Item(9, 17).Value = 41.85674589533234
$ws.Cells.Item(9, 18).Value = 376.710713057991
$ws.Cells.Item(9, 19).Value = 0.005666267645992037
$ws.Cells.Item(9, 20).Value = 0.005666267645992034
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.041671
$ws.Cells.Item(10, 8).Value = 9.125013000000001
$ws.Cells.Item(10, 9).Value = 0.03538979513654388
$ws.Cells.Item(10, 10).Value = 0.03538979513654388
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.041671
$ws.Cells.Item(10, 14).Value = 9.125013000000001
$ws.Cells.Item(10, 15).Value = 0.03538979513654388
$ws.Cells.Item(10, 16).Value = 0.03538979513654388
$ws.Cells.Item(10, 17).Value = 9.251762472241003
$ws.Cells.Item(10, 18).Value = 9.251762472241003
$ws.Cells.Item(10, 19).Value = 0.001252437599806545
$ws.Cells.Item(10, 20).Value = 0.001252437599806545
